# Update Excel files after daily scrape - 2025-11-11 03:24:35 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths ---
# NOTE: this runtime's ColumnWidth -> stored <col width> mapping adds a
# constant 5/6 padding offset (width_attr = ColumnWidth + 5/6), so we
# subtract it here to land on the exact target width stored in the XML.
$widthPad = 5 / 6
$ws.Columns.Item(3).ColumnWidth = 69 - $widthPad
$ws.Columns.Item(4).ColumnWidth = 30 - $widthPad
$ws.Columns.Item(7).ColumnWidth = 16 - $widthPad
$ws.Columns.Item(8).ColumnWidth = 41 - $widthPad

# --- Row data (A..H) for rows 2-7 ---
$rows = @(
    @{ Row = 2;  A = "1329108"; B = "https://aiesec.org/opportunity/global-talent/1329108"; C = "[Partly Remote] Occupational Health and Safety Projects Specialist"; D = "Mexico City, CDMX, Mexico";   E = "No"; F = "0 applicants";  G = "Partly Remote";  H = "Sodexo Mexico" },
    @{ Row = 3;  A = "1329104"; B = "https://aiesec.org/opportunity/global-talent/1329104"; C = "Sales Intern";                                                               D = "Pannipitiya, Sri Lanka";  E = "No"; F = "3 applicants";  G = "6 - 18 Months";  H = "Frella International" },
    @{ Row = 4;  A = "1328030"; B = "https://aiesec.org/opportunity/global-talent/1328030"; C = "[Remote] Software Application Support and Development";                      D = "No location available";   E = "No"; F = "17 applicants"; G = "Remote";         H = "dJava Factory Sdn Bhd" },
    @{ Row = 5;  A = "1328026"; B = "https://aiesec.org/opportunity/global-talent/1328026"; C = "[Remote] Software Application Support and Development";                      D = "No location available";   E = "No"; F = "32 applicants"; G = "Remote";         H = "dJava Factory Sdn Bhd" },
    @{ Row = 6;  A = "1327354"; B = "https://aiesec.org/opportunity/global-talent/1327354"; C = "Guest Relations Executive";                                                   D = "Weligama, Sri Lanka";     E = "No"; F = "13 applicants"; G = "3 - 6 Months";   H = "Steradian Capital Investments" },
    @{ Row = 7;  A = "1327154"; B = "https://aiesec.org/opportunity/global-talent/1327154"; C = "Guest Relations and Service Captain";                                         D = "Hyderabad, Telangana, India"; E = "No"; F = "8 applicants"; G = "3 - 6 Months"; H = "Concu [Quenelle36 Foodworks Pvt. Ltd.]" }
)

# Column A holds numeric-looking opportunity IDs that must stay text
# (matching the source data, which is plain text, not a number). Assigning
# a numeric-looking string straight to .Value auto-converts it to a real
# number, so mark the cell as Text first, then reset the style back to
# Normal so no stray style index is left attached to the cell.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

foreach ($r in $rows) {
    $rowIndex = $r.Row
    Set-TextValue $rowIndex 1 $r.A
    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
    $ws.Cells.Item($rowIndex, 5).Value = $r.E
    $ws.Cells.Item($rowIndex, 6).Value = $r.F
    $ws.Cells.Item($rowIndex, 7).Value = $r.G
    $ws.Cells.Item($rowIndex, 8).Value = $r.H
}
